# Update Consumption_Actual_Historical workbook:
#  - shift all timestamps (column A) forward by 28 days (20/21 Aug 2025 -> 17/18 Sep 2025)
#  - refresh the "Actual Consumption (MW)" values (column B) with the newly retrained figures
#  - the "Lookup" text in column D is derived from the date, so it updates automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Consumption (MW) values for rows 2..193 (row 1 is the header)
$bValues = @(
    5273,5270,5185,5150,5145,5162,5126,5097,5010,5055,
    5110,5088,5120,5077,5074,5121,5166,5267,5301,5337,
    5515,5659,5763,5887,6085,6265,6421,6503,6640,6770,
    6796,6739,6802,6957,6934,6836,6676,6655,6627,6505,
    6372,6279,6254,6169,5998,5966,5956,5995,6002,5947,
    5955,5842,5847,5804,5725,5681,5765,5768,5768,5723,
    5661,5671,5764,5805,5861,5954,6015,6077,6240,6361,
    6498,6524,6572,6729,6904,7021,7059,7063,7048,7020,
    6901,6806,6698,6574,6365,6292,6183,6049,5773,5627,
    5566,5490,5371,5275,5200,5184,5184,5139,5137,5095,
    4961,5014,4980,4962,4860,4926,4970,4926,4970,4936,
    5016,4971,5096,5138,5151,5247,5446,5551,5678,5764,
    5974,6151,6227,6263,6255,6172,6178,6171,5883,5818,
    5783,5616,5324,5195,5135,5055,4914,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0
)

$lastRow = 193
$dayShift = 28

for ($row = 2; $row -le $lastRow; $row++) {
    $aCell = $ws.Cells.Item($row, 1)
    $newDate = $aCell.Value2 + $dayShift
    $aCell.Value2 = $newDate

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value2 = $bValues[$row - 2]

    $cValue = $ws.Cells.Item($row, 3).Value2
    $dateText = $excel.WorksheetFunction.Text($newDate, "dd.mm.yyyy")
    $ws.Cells.Item($row, 4).Value = $dateText + $cValue
}
